# 🔄 MAJ automatique BRVM via GitHub Actions
#
# Refreshes the "Recommandations" ranking table (re-sorted by Variation
# Totale (%) descending, as produced by the scraper) and the companion
# "Top_YTD" leaderboard on every scheduled run.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Freshly computed "Recommandations" rows (A:G), in their new sorted order.
# Columns: Titre, Jours en Hausse, Jours en Baisse, Variation Totale (%),
#          Derniere Variation (%), Recommandation, Strategie
$recommandations = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3378.18, 115.5, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2810, 700, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2304.25, 596.33, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1979.33, 496.01, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI", 0, 2, 1850, 925, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1644.86, 416.11, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1333.53, 338.4, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 1, 940, 940, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 704.4, 176.33, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 575.77, 143.57, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 565.87, 141.1, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 561.84, 140.62, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 516.98, 134.35, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 426.83, 107.64, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 379.22, 94.93, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                 (**)", 0, 1, 230.07, 230.07, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                  (**)", 0, 1, 226.03, 226.03, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                     (**)", 0, 1, 209.18, 209.18, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 1, 208.32, 208.32, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE            (**)", 0, 1, 197.25, 197.25, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE          (**)", 0, 1, 193.92, 193.92, "🟡 Observer", "➖ Neutre"),
    @("FILTISAC CI (FTSC)", 4, 0, 29.55, 7.4, "🟢 Achat", "✅ Renforcer"),
    @("UNIWAX CI (UNXC)", 3, 1, 16.83, -5.06, "🟢 Achat", "✅ Renforcer"),
    @("PALM CI (PALC)", 2, 0, 14.77, 7.47, "🟡 Observer", "➖ Neutre"),
    @("SICABLE CI (CABC)", 2, 0, 14.65, 7.32, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 1, 0, 7.49, 7.49, "🟡 Observer", "➖ Neutre"),
    @("VIVO ENERGY CI (SHEC)", 1, 0, 7.34, 7.34, "🟡 Observer", "➖ Neutre"),
    @("SERVAIR ABIDJAN CI (ABJC)", 1, 0, 6.6, 6.6, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 2, 2, 0.15, 7.38, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 2, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("NESTLE CI (NTLC)", 1, 1, -0.15, 7.32, "🟡 Observer", "👀 À surveiller"),
    @("SOLIBRA CI (SLBC)", 0, 1, -1.2, -1.2, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 1, 1, -1.63, 5.41, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA BN (BOAB)", 0, 1, -1.92, -1.92, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 0, 1, -2.73, -2.73, "🟡 Observer", "➖ Neutre"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -3.08, -3.08, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 0, 1, -4.17, -4.17, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 1, 2, -7.35, -7.39, "🟡 Observer", "👀 À surveiller"),
    @("SICOR CI (SICC)", 0, 1, -7.38, -7.38, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 1, 3, -14.82, 7.19, "🔴 Vente", "⚠️ Risque de décrochage"),
    @("BERNABE CI (BNBC)", 0, 3, -19.69, -4.91, "🔴 Vente", "⚠️ Risque de décrochage")
)

# Freshly computed "Top_YTD" rows (A:B): Titre, Progression YTD (%)
$topYtd = @(
    @("BRVM - SERVICES PUBLICS", 10145701.29),
    @("AIR LIQUIDE CI", 414572),
    @("BRVM - AUTRES SECTEURS", 208600.98),
    @("BRVM - DISTRIBUTION", 125032.2),
    @("BRVM - TRANSPORT", 68186.48),
    @("BRVM - AGRICULTURE", 35153.57),
    @("SETAO CI", 10406.25),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 5707.24),
    @("BRVM - FINANCES", 3441.17),
    @("BRVM - SERVICES FINANCIERS", 3299.63)
)

$startRow = 2
for ($i = 0; $i -lt $recommandations.Count; $i++) {
    $row = $startRow + $i
    $values = $recommandations[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws1.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# The refreshed table only has 41 data rows (42 incl. header) this run,
# two fewer than before (44) — clear the now-stale trailing rows so the
# used range shrinks to A1:G42 instead of leaving orphaned values behind.
$lastOldRow = $startRow + 42
$ws1.Range("A" + ($startRow + $recommandations.Count) + ":G" + $lastOldRow).ClearContents()

for ($i = 0; $i -lt $topYtd.Count; $i++) {
    $row = $startRow + $i
    $values = $topYtd[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws2.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
